# Update crypto price/volume figures per the Thu Sep 19 21:54:30 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.974.81"
$ws.Range("E2").Value = "  +4.43%  "

$ws.Range("D3").Value = "2.464.91"
$ws.Range("E3").Value = "  +5.58%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'564.33"

$ws.Range("D6").Value = "'142.57"
$ws.Range("E6").Value = "  +8.62%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.588"
$ws.Range("E8").Value = "  +1.36%  "

$ws.Range("D9").Value = "2.463.98"
$ws.Range("E9").Value = "  +5.61%  "

$ws.Range("E10").Value = "  +3.38%  "

$ws.Range("E11").Value = "  +1.74%  "

$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("D13").Value = "'0.351"
$ws.Range("E13").Value = "  +4.27%  "

$ws.Range("D14").Value = "'26.25"
$ws.Range("E14").Value = "  +10.69%  "

$ws.Range("D15").Value = "2.905.37"
$ws.Range("E15").Value = "  +5.64%  "

$ws.Range("D16").Value = "62.885.95"
$ws.Range("E16").Value = "  +4.33%  "

$ws.Range("E17").Value = "  +4.56%  "

$ws.Range("D18").Value = "2.465.69"
$ws.Range("E18").Value = "  +5.75%  "

$ws.Range("E19").Value = "  +5.22%  "

$ws.Range("D20").Value = "'340.42"
$ws.Range("E20").Value = "  +8.19%  "

$ws.Range("E21").Value = "  +3.83%  "

$ws.Range("E22").Value = "  +3.05%  "

$ws.Range("D23").Value = "'0.999"

$ws.Range("D24").Value = "'65.49"
$ws.Range("E24").Value = "  +2.09%  "

$ws.Range("E25").Value = "  +1.44%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").Value = "'1.50"
$ws.Range("E27").Value = "  +7.93%  "

$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("D29").Value = "'1.37"
$ws.Range("E29").Value = "  +8.58%  "

$ws.Range("E30").Value = "  +11.63%  "

$ws.Range("E31").Value = "  +9.34%  "

$ws.Range("E32").Value = "  +6.23%  "

$ws.Range("D33").Value = "'176.37"
$ws.Range("E33").Value = "  +2.98%  "

$ws.Range("D34").Value = "'1.51"
$ws.Range("E34").Value = "  +10.58%  "

$ws.Range("D35").Value = "'0.396"
$ws.Range("E35").Value = "  +2.96%  "

$ws.Range("D36").Value = "'18.83"
$ws.Range("E36").Value = "  +3.98%  "

$ws.Range("D37").Value = "'365.36"
$ws.Range("E37").Value = "  +12.72%  "

$ws.Range("E38").Value = "  +7.12%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("E41").Value = "  +10.31%  "

$ws.Range("D42").Value = "'40.46"
$ws.Range("E42").Value = "  +6.14%  "

$ws.Range("D43").Value = "'149.67"
$ws.Range("E43").Value = "  +8.67%  "

$ws.Range("E44").Value = "  +5.29%  "

$ws.Range("D45").Value = "'20.54"
$ws.Range("E45").Value = "  +6.27%  "

$ws.Range("D46").Value = "'0.596"
$ws.Range("E46").Value = "  +5.44%  "

$ws.Range("E47").Value = "  +0.69%  "

$ws.Range("E48").Value = "  +3.24%  "

$ws.Range("D49").Value = "0.0₆0239"
$ws.Range("E49").Value = "  +9.33%  "

$ws.Range("E50").Value = "  +4.41%  "

$ws.Range("D51").Value = "'17.93"
$ws.Range("E51").Value = "  +4.73%  "

